# "Add column value 33% + Clean Code"
#
# The schedule table (columns A:C, "Desenvolvedor"/"Inicio"/"fim") currently
# ends at row 197. We extend it with 6 more rows (198-203, ~33% more rows)
# and make sure the "end of table" styling (the bottom border band that
# currently decorates row 197) moves down to the new final row (203)
# instead, so row 197 becomes a normal interior data row.
#
# This runtime's Range.Copy(destination) form (copy with an explicit
# destination argument) is the form that reliably carries the source
# range's cell-level formatting (number format, borders, fill, alignment)
# to the destination, so it's used here to stamp the right look onto every
# new/re-styled row before the actual text values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot today's "last row" look (the bottom-border band + highlight
#    fill currently on row 197) down onto what will become the new last
#    row, 203, BEFORE row 197 itself is restyled as a normal row.
$ws.Range("A197:C197").Copy($ws.Range("A203:C203"))

# 2) Row 197 turns into an ordinary data row, matching the style already
#    used one row above it (row 196).
$ws.Range("A196:C196").Copy($ws.Range("A197:C197"))

# 3) Row 198 (first brand-new row) uses that same ordinary-row style.
$ws.Range("A196:C196").Copy($ws.Range("A198:C198"))

# 4) Rows 199-202 use the alternate interior style already used by row 191
#    (and similar rows throughout the table).
$ws.Range("A191:C191").Copy($ws.Range("A199:C199"))
$ws.Range("A191:C191").Copy($ws.Range("A200:C200"))
$ws.Range("A191:C191").Copy($ws.Range("A201:C201"))
$ws.Range("A191:C191").Copy($ws.Range("A202:C202"))

# All data rows in this table share the same custom row height.
$ws.Rows.Item(198).RowHeight = 19.55
$ws.Rows.Item(199).RowHeight = 19.55
$ws.Rows.Item(200).RowHeight = 19.55
$ws.Rows.Item(201).RowHeight = 19.55
$ws.Rows.Item(202).RowHeight = 19.55
$ws.Rows.Item(203).RowHeight = 19.55

# Now fill in the actual schedule data for rows 197-203.
# Row 197 keeps its original data (Matheus 11:00 -> 17:00), only the look
# changed above.
$ws.Range("A197").Value = "Matheus"
$ws.Range("B197").Value = "11:00:00"
$ws.Range("C197").Value = "17:00:00"

# Row 198: new developer, Ana.
$ws.Range("A198").Value = "Ana"
$ws.Range("B198").Value = "17:00:00"
$ws.Range("C198").Value = "22:30:00"

# Rows 199-202: additional evening-shift entries.
$ws.Range("A199").Value = "Jonas"
$ws.Range("B199").Value = "18:00:00"
$ws.Range("C199").Value = "22:30:00"

$ws.Range("A200").Value = "Matheus"
$ws.Range("B200").Value = "18:00:00"
$ws.Range("C200").Value = "22:30:00"

$ws.Range("A201").Value = "Thiago"
$ws.Range("B201").Value = "18:00:00"
$ws.Range("C201").Value = "22:30:00"

$ws.Range("A202").Value = "Matheus"
$ws.Range("B202").Value = "18:00:00"
$ws.Range("C202").Value = "22:30:00"

# Row 203: new final row of the table.
$ws.Range("A203").Value = "Thiago"
$ws.Range("B203").Value = "18:00:00"
$ws.Range("C203").Value = "22:30:00"

Write-Host "Added rows 198-203 (incl. new developer 'Ana'); table now spans A1:C203 with row 203 as the new styled last row."
